# Add three new account rows to the "Export" sheet.
# New rows (Conta, Nome, Saldo):
#   004482163, CIRIA,      25280.81   -> inserted right above the MURYLO row
#   004948033, GUILHERME,  25000      -> inserted right above the MURYLO row
#   005002457, ROSANGELA,  2778.74    -> inserted right after the RODRIGO/5000 row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert CIRIA / 25280.81 at row 2 (pushes MURYLO and everything below down) ---
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004482163"
$ws.Cells.Item(2, 2).Value = "CIRIA"
$ws.Cells.Item(2, 3).Value = 25280.81

# --- Insert GUILHERME / 25000 at row 3 (directly below CIRIA, above MURYLO) ---
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004948033"
$ws.Cells.Item(3, 2).Value = "GUILHERME"
$ws.Cells.Item(3, 3).Value = 25000

# --- Insert ROSANGELA / 2778.74 at row 6, right after the RODRIGO (5000) row ---
# (Row order now: 1 header, 2 CIRIA, 3 GUILHERME, 4 MURYLO, 5 RODRIGO/5000, 6 <new>)
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "005002457"
$ws.Cells.Item(6, 2).Value = "ROSANGELA"
$ws.Cells.Item(6, 3).Value = 2778.74
